$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so numeric-looking strings
# (e.g. "245.13", "-0.68%") are preserved verbatim as text, matching the source data.
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("D2").Value = '245.13'
$ws.Range("E2").Value = '-0.68%'
$ws.Range("D3").Value = '29.08'
$ws.Range("E3").Value = '-1.85%'
$ws.Range("D4").Value = '5.273'
$ws.Range("E4").Value = '1.57%'
$ws.Range("D5").Value = '0.05714'
$ws.Range("E5").Value = '0.08%'
$ws.Range("D6").Value = '6.605'
$ws.Range("E6").Value = '0.37%'
$ws.Range("D7").Value = '3.175'
$ws.Range("E7").Value = '3.49%'
$ws.Range("D8").Value = '0.8532'
$ws.Range("E8").Value = '-0.63%'
$ws.Range("D9").Value = '0.8580'
$ws.Range("E9").Value = '-2.36%'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '0.01009'
$ws.Range("E10").Value = '1,577.82%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1373'
$ws.Range("E11").Value = '0.25%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.07072'
$ws.Range("E12").Value = '0.00%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03166'
$ws.Range("E13").Value = '10.04%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09289'
$ws.Range("E14").Value = '-1.02%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001535'
$ws.Range("E15").Value = '1.44%'
$ws.Range("D16").Value = '0.005934'
$ws.Range("E16").Value = '-4.14%'
$ws.Range("D17").Value = '3.508'
$ws.Range("E17").Value = '0.83%'
$ws.Range("E18").Value = '-4.40%'
$ws.Range("D19").Value = '0.3172'
$ws.Range("E19").Value = '-0.14%'
$ws.Range("D20").Value = '0.03322'
$ws.Range("E20").Value = '0.48%'
$ws.Range("E21").Value = '-1.73%'
$ws.Range("D22").Value = '3.486'
$ws.Range("E22").Value = '0.46%'
$ws.Range("D23").Value = '0.04120'
$ws.Range("E23").Value = '-1.44%'
$ws.Range("D24").Value = '0.1328'
$ws.Range("E24").Value = '-3.62%'
$ws.Range("D25").Value = '0.001218'
$ws.Range("E25").Value = '0.00%'
$ws.Range("D26").Value = '0.004146'
$ws.Range("E26").Value = '-17.92%'
$ws.Range("E27").Value = '-0.78%'
$ws.Range("D28").Value = '0.0001448'
$ws.Range("E28").Value = '-25.29%'
$ws.Range("D40").Value = '0.03792'
$ws.Range("E40").Value = '0.94%'
$ws.Range("D41").Value = '0.1063'
$ws.Range("E41").Value = '-0.95%'
$ws.Range("D42").Value = '0.002415'
$ws.Range("E42").Value = '15.11%'
$ws.Range("E43").Value = '-48.04%'
$ws.Range("D44").Value = '0.009389'
$ws.Range("E44").Value = '-5.61%'
$ws.Range("D45").Value = '0.00005263'
$ws.Range("E45").Value = '2.65%'
$ws.Range("E46").Value = '0.04%'
$ws.Range("D47").Value = '0.08985'
$ws.Range("E47").Value = '26.65%'
$ws.Range("D48").Value = '0.002436'
$ws.Range("E48").Value = '-5.69%'
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").Value = '0.04%'
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").Value = '0.04%'
